$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C4").Value = "0h 59m"
$ws.Range("A4").Value = "web design"
$ws.Range("A5").Select()
